$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2991.8857
$ws.Range("I64").Value = 2784.5454
$ws.Range("J64").Value = 3342.7693
$ws.Range("K64").Value = 2784.5454
$ws.Range("L64").Value = 3342.7693
$ws.Range("M64").Value = -2536.5454
$ws.Range("N64").Value = -3838.7693
$ws.Range("H67").Value = 2991.8857
$ws.Range("I67").Value = 2784.5454
$ws.Range("J67").Value = 3342.7693
$ws.Range("K67").Value = 2784.5454
$ws.Range("L67").Value = 3342.7693
$ws.Range("M67").Value = -1926.5454
$ws.Range("N67").Value = -5058.7693
$ws.Range("H76").Value = 34485364
$ws.Range("I76").Value = 45457124
$ws.Range("J76").Value = 2700
$ws.Range("K76").Value = 45457124
$ws.Range("L76").Value = 2700
$ws.Range("M76").Value = -45456809
$ws.Range("N76").Value = -3330
$ws.Range("H79").Value = 34485364
$ws.Range("I79").Value = 45457124
$ws.Range("J79").Value = 2700
$ws.Range("K79").Value = 45457124
$ws.Range("L79").Value = 2700
$ws.Range("M79").Value = -45456032
$ws.Range("N79").Value = -4884
$ws.Range("H82").Value = 4858.3335
$ws.Range("I82").Value = 383.33334
$ws.Range("J82").Value = 9333.333000000001
$ws.Range("K82").Value = 1150.00002
$ws.Range("L82").Value = 27999.999
$ws.Range("M82").Value = -744.0000199999999
$ws.Range("N82").Value = -28811.999
$ws.Range("H85").Value = 4858.3335
$ws.Range("I85").Value = 383.33334
$ws.Range("J85").Value = 9333.333000000001
$ws.Range("K85").Value = 1150.00002
$ws.Range("L85").Value = 27999.999
$ws.Range("M85").Value = 253.9999800000001
$ws.Range("N85").Value = -30807.999
$ws.Range("H106").Value = 55556350
$ws.Range("I106").Value = 58824336
$ws.Range("J106").Value = 600
$ws.Range("K106").Value = 58824336
$ws.Range("L106").Value = 600
$ws.Range("M106").Value = -58823705
$ws.Range("N106").Value = -1862
$ws.Range("H137").Value = 16282578
$ws.Range("I137").Value = 3907108.8
$ws.Range("J137").Value = 55884080
$ws.Range("K137").Value = 11721326.4
$ws.Range("L137").Value = 167652240
$ws.Range("M137").Value = -11718776.4
$ws.Range("N137").Value = -167657340

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5396.34
$ws.Range("I32").Value = 4307.8926
$ws.Range("J32").Value = 19857.143
$ws.Range("K32").Value = 4307.8926
$ws.Range("L32").Value = 19857.143
$ws.Range("M32").Value = -4020.8926
$ws.Range("N32").Value = -20431.143
$ws.Range("H63").Value = 1948.0435
$ws.Range("I63").Value = 1919.2858
$ws.Range("K63").Value = 1919.2858
$ws.Range("M63").Value = -1233.2858
$ws.Range("H66").Value = 1948.0435
$ws.Range("I66").Value = 1919.2858
$ws.Range("K66").Value = 9596.429
$ws.Range("M66").Value = -6164.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1731.3334
$ws.Range("I105").Value = 1765.5555
$ws.Range("J105").Value = 1680
$ws.Range("K105").Value = 1765.5555
$ws.Range("L105").Value = 1680
$ws.Range("M105").Value = -18.55549999999994
$ws.Range("N105").Value = -5174

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 13891488
$ws.Range("I62").Value = 2252.7273
$ws.Range("J62").Value = 35717428
$ws.Range("K62").Value = 2252.7273
$ws.Range("L62").Value = 35717428
$ws.Range("M62").Value = -1628.7273
$ws.Range("N62").Value = -35718676
$ws.Range("H65").Value = 13891488
$ws.Range("I65").Value = 2252.7273
$ws.Range("J65").Value = 35717428
$ws.Range("K65").Value = 11263.6365
$ws.Range("L65").Value = 178587140
$ws.Range("M65").Value = -8143.636500000001
$ws.Range("N65").Value = -178593380
$ws.Range("H99").Value = 9868.714
$ws.Range("I99").Value = 16264.8
$ws.Range("J99").Value = 8478.261
$ws.Range("K99").Value = 16264.8
$ws.Range("L99").Value = 8478.261
$ws.Range("M99").Value = -14766.8
$ws.Range("N99").Value = -11474.261
$ws.Range("H126").Value = 9868.714
$ws.Range("I126").Value = 16264.8
$ws.Range("J126").Value = 8478.261
$ws.Range("K126").Value = 48794.39999999999
$ws.Range("L126").Value = 25434.783
$ws.Range("M126").Value = -46324.39999999999
$ws.Range("N126").Value = -30374.783

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 3582.0527
$ws.Range("I75").Value = 1428.25
$ws.Range("J75").Value = 4156.4
$ws.Range("K75").Value = 4284.75
$ws.Range("L75").Value = 12469.2
$ws.Range("M75").Value = -3286.75
$ws.Range("N75").Value = -14465.2
$ws.Range("H78").Value = 3582.0527
$ws.Range("I78").Value = 1428.25
$ws.Range("J78").Value = 4156.4
$ws.Range("K78").Value = 12854.25
$ws.Range("L78").Value = 37407.6
$ws.Range("M78").Value = -7862.25
$ws.Range("N78").Value = -47391.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5708458
$ws.Range("I70").Value = 2845265.8
$ws.Range("J70").Value = 10207760
$ws.Range("K70").Value = 2845265.8
$ws.Range("L70").Value = 10207760
$ws.Range("M70").Value = -2844995.8
$ws.Range("N70").Value = -10208300
$ws.Range("H73").Value = 5708458
$ws.Range("I73").Value = 2845265.8
$ws.Range("J73").Value = 10207760
$ws.Range("K73").Value = 2845265.8
$ws.Range("L73").Value = 10207760
$ws.Range("M73").Value = -2844329.8
$ws.Range("N73").Value = -10209632
$ws.Range("H80").Value = 7746.4707
$ws.Range("I80").Value = 3924.1667
$ws.Range("J80").Value = 16920
$ws.Range("K80").Value = 3924.1667
$ws.Range("L80").Value = 16920
$ws.Range("M80").Value = -2926.1667
$ws.Range("N80").Value = -18916
$ws.Range("H83").Value = 7746.4707
$ws.Range("I83").Value = 3924.1667
$ws.Range("J83").Value = 16920
$ws.Range("K83").Value = 19620.8335
$ws.Range("L83").Value = 84600
$ws.Range("M83").Value = -14628.8335
$ws.Range("N83").Value = -94584

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 41669620
$ws.Range("I46").Value = 2119.75
$ws.Range("J46").Value = 62503370
$ws.Range("K46").Value = 2119.75
$ws.Range("L46").Value = 62503370
$ws.Range("M46").Value = -1931.75
$ws.Range("N46").Value = -62503746
$ws.Range("H50").Value = 20616.666
$ws.Range("J50").Value = 20616.666
$ws.Range("L50").Value = 20616.666
$ws.Range("N50").Value = -21890.666
$ws.Range("H68").Value = 1829.3158
$ws.Range("I68").Value = 1172
$ws.Range("K68").Value = 1172
$ws.Range("M68").Value = -423
$ws.Range("H71").Value = 1829.3158
$ws.Range("I71").Value = 1172
$ws.Range("K71").Value = 5860
$ws.Range("M71").Value = -2116

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 93767580
$ws.Range("J62").Value = 6640
$ws.Range("L62").Value = 6640
$ws.Range("N62").Value = -7888
$ws.Range("H65").Value = 93767580
$ws.Range("J65").Value = 6640
$ws.Range("L65").Value = 33200
$ws.Range("N65").Value = -39440
